$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..26, columns A..Q (A = model name string, B..Q = numeric metrics)
$data = @(
    ,("model_8_2_0", 0.9999989699042762, 0.9990244677798398, 0.9999684164407366, 0.9999987730294633, 0.9999958254641106, 0.0000009615491421736449, 0.0009106165066041113, 0.000005445466703584143, 0.0000005642146659887532, 0.000003004840684786449, 0.00005129893614687497, 0.0009805861217525185, 1.000024722297372, 0.001022331703048629, 77.70944032663732, 108.1813359483423)
    ,("model_8_2_22", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_21", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_20", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_19", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_18", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_17", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_16", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_15", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_14", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_13", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_23", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_12", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_10", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_9", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_8", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_7", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_6", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_5", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_4", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_3", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_2", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_1", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_11", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
    ,("model_8_2_24", 0.9999989698596451, 0.9990244665669379, 0.9999684164407366, 0.9999987724631316, 0.9999958252832108, 0.0000009615908033740045, 0.0009106176387947784, 0.000005445466703584143, 0.0000005644750900211889, 0.000003004970896802666, 0.00005130350677945762, 0.0009806073645318011, 1.000024723368518, 0.001022353850177004, 77.70935367418031, 108.1812492958853)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    for ($c = 1; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
